# Add I0 and IF columns (I and J) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold/centered/bordered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-29
$data = @(
    @(2, 7, 8),
    @(3, 6, 7),
    @(4, 7, 8),
    @(5, 7, 8),
    @(6, 8, 8),
    @(7, 7, 8),
    @(8, 8, 8),
    @(9, 9, 9),
    @(10, 8, 9),
    @(11, 7, 8),
    @(12, 7, 8),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 8, 8),
    @(17, 5, 6),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 6, 6),
    @(21, 9, 9),
    @(22, 7, 7),
    @(23, 9, 9),
    @(24, 5, 5),
    @(25, 1, 4),
    @(26, 2, 4),
    @(27, 4, 4),
    @(28, 1, 4),
    @(29, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
